$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 55.687
$ws.Range("D2").Value = 55.687
$ws.Range("E2").Value = 3.08659142
$ws.Range("F2").Value = 0.01866754
$ws.Range("G2").Value = 1.02818196
$ws.Range("H2").Value = 57.55636339
$ws.Range("I2").Value = 7.186980150500186
$ws.Range("J2").Value = 7.186980150500186
$ws.Range("K2").Value = 0.4096919970236285
$ws.Range("L2").Value = 0.002414330495218113
$ws.Range("M2").Value = 0.1115499255062003
$ws.Range("N2").Value = 11.15901585739432

# Row 3
$ws.Range("C3").Value = 89.771
$ws.Range("D3").Value = 89.771
$ws.Range("E3").Value = 1.91476307
$ws.Range("F3").Value = 0.00684953
$ws.Range("G3").Value = 0.60273417
$ws.Range("H3").Value = 53.93090927
$ws.Range("I3").Value = 11.80660641145266
$ws.Range("J3").Value = 11.80660641145266
$ws.Range("K3").Value = 0.2492017269679303
$ws.Range("L3").Value = 0.001380181430059063
$ws.Range("M3").Value = 0.08037549784832107
$ws.Range("N3").Value = 9.09044961075559

# Row 4
$ws.Range("C4").Value = 27.741
$ws.Range("D4").Value = 55.457
$ws.Range("E4").Value = 3.18733622
$ws.Range("F4").Value = 0.03378314
$ws.Range("G4").Value = 0.45638717
$ws.Range("H4").Value = 12.8867726
$ws.Range("I4").Value = 5.941146035590794
$ws.Range("J4").Value = 11.88020581851553
$ws.Range("K4").Value = 0.6784729570128122
$ws.Range("L4").Value = 0.005410311859293771
$ws.Range("M4").Value = 0.06156511207588285
$ws.Range("N4").Value = 4.062709615798074

# Row 5
$ws.Range("C5").Value = 48.033
$ws.Range("D5").Value = 94.18
$ws.Range("E5").Value = 1.8425729
$ws.Range("F5").Value = 0.0109313
$ws.Range("G5").Value = 0.25438118
$ws.Range("H5").Value = 12.18384207
$ws.Range("I5").Value = 8.351626228505479
$ws.Range("J5").Value = 15.48933567069678
$ws.Range("K5").Value = 0.3016741180042857
$ws.Range("L5").Value = 0.002570304373916555
$ws.Range("M5").Value = 0.03720237371213586
$ws.Range("N5").Value = 2.590364271959498

# Row 6
$ws.Range("C6").Value = 14.436
$ws.Range("D6").Value = 57.699
$ws.Range("E6").Value = 3.13267109
$ws.Range("F6").Value = 0.05759362
$ws.Range("G6").Value = 0.20001618
$ws.Range("H6").Value = 2.96714329
$ws.Range("I6").Value = 3.768568177186941
$ws.Range("J6").Value = 15.07360648990853
$ws.Range("K6").Value = 0.8157236465990996
$ws.Range("L6").Value = 0.01089640248351183
$ws.Range("M6").Value = 0.0329796359439366
$ws.Range("N6").Value = 1.153716206375176

# Row 7
$ws.Range("B7").Value = 0.9999777799999999
$ws.Range("C7").Value = 25.443
$ws.Range("D7").Value = 93.635
$ws.Range("E7").Value = 1.86505701
$ws.Range("F7").Value = 0.01612009
$ws.Range("G7").Value = 0.09794596000000001
$ws.Range("H7").Value = 2.49143973
$ws.Range("I7").Value = 5.731981036832433
$ws.Range("J7").Value = 17.13883060522259
$ws.Range("K7").Value = 0.3389970062386963
$ws.Range("L7").Value = 0.004702166826556111
$ws.Range("M7").Value = 0.01990200645008284
$ws.Range("N7").Value = 0.7683789732638323

# Row 8
$ws.Range("C8").Value = 9.359
$ws.Range("D8").Value = 56.07
$ws.Range("E8").Value = 3.24010606
$ws.Range("F8").Value = 0.08475283
$ws.Range("G8").Value = 0.12835513
$ws.Range("H8").Value = 1.25425688
$ws.Range("I8").Value = 2.686137437722488
$ws.Range("J8").Value = 16.10150174233318
$ws.Range("K8").Value = 0.8534073801664708
$ws.Range("L8").Value = 0.01432791601274827
$ws.Range("M8").Value = 0.02789220441001302
$ws.Range("N8").Value = 0.6122301681147503

# Row 9
$ws.Range("C9").Value = 17.421
$ws.Range("D9").Value = 88.746
$ws.Range("E9").Value = 1.96655922
$ws.Range("F9").Value = 0.01782611
$ws.Range("G9").Value = 0.049275
$ws.Range("H9").Value = 0.86699929
$ws.Range("I9").Value = 4.495327236916904
$ws.Range("J9").Value = 16.17352982706827
$ws.Range("K9").Value = 0.3529428848132966
$ws.Range("L9").Value = 0.005597369874987485
$ws.Range("M9").Value = 0.01204148100795184
$ws.Range("N9").Value = 0.3718073392068585

# Row 10
$ws.Range("C10").Value = 7.042
$ws.Range("D10").Value = 56.2
$ws.Range("E10").Value = 3.19587531
$ws.Range("F10").Value = 0.10072949
$ws.Range("G10").Value = 0.08720381999999999
$ws.Range("H10").Value = 0.6403016699999999
$ws.Range("I10").Value = 1.77859478223664
$ws.Range("J10").Value = 14.19131401221865
$ws.Range("K10").Value = 0.7936199990741376
$ws.Range("L10").Value = 0.01376534630144475
$ws.Range("M10").Value = 0.01922850157137265
$ws.Range("N10").Value = 0.2870390857682448

# Row 11
$ws.Range("C11").Value = 13.122
$ws.Range("D11").Value = 80.984
$ws.Range("E11").Value = 2.15736031
$ws.Range("F11").Value = 0.01699851
$ws.Range("G11").Value = 0.02655863
$ws.Range("H11").Value = 0.35377627
$ws.Range("I11").Value = 3.395385787306745
$ws.Range("J11").Value = 14.81780386262793
$ws.Range("K11").Value = 0.39833869430396
$ws.Range("L11").Value = 0.006140531899755258
$ws.Range("M11").Value = 0.008355805724908426
$ws.Range("N11").Value = 0.1625674391459051

# Row 12
$ws.Range("C12").Value = 5.703
$ws.Range("D12").Value = 56.846
$ws.Range("E12").Value = 3.18929553
$ws.Range("F12").Value = 0.1163153
$ws.Range("G12").Value = 0.06571637000000001
$ws.Range("H12").Value = 0.39689023
$ws.Range("I12").Value = 1.532690490653654
$ws.Range("J12").Value = 15.29209377855538
$ws.Range("K12").Value = 0.8581139051045921
$ws.Range("L12").Value = 0.01456758793102441
$ws.Range("M12").Value = 0.01706750892145318
$ws.Range("N12").Value = 0.2031138917855777

# Row 13
$ws.Range("C13").Value = 10.654
$ws.Range("D13").Value = 74.155
$ws.Range("E13").Value = 2.35904075
$ws.Range("F13").Value = 0.01551286
$ws.Range("G13").Value = 0.01572982
$ws.Range("H13").Value = 0.17231472
$ws.Range("I13").Value = 3.126989296761673
$ws.Range("J13").Value = 14.04607657914719
$ws.Range("K13").Value = 0.4366621609383309
$ws.Range("L13").Value = 0.006359856239868303
$ws.Range("M13").Value = 0.006103465512219993
$ws.Range("N13").Value = 0.09866405602686965
